$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row
$ws.Range("C1").Value = "JENIS KELAMIN"
$ws.Range("D1").Value = "ALAMAT"

# Update data row
$ws.Range("A2").Value = 12
$ws.Range("B2").Value = "Puput"
$ws.Range("C2").Value = "Perempuan"
$ws.Range("D2").Value = "Pengkok"
